# Auto-generated Word COM-interop script implementing the diff.
$d = $word.ActiveDocument

# Change1 - formulario
$rng = $d.Content
$found = $rng.Find.Execute('El sistema tendrá un formulario el cual tiene que llenar para ir agregando los productos al histórico diario.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  Write-Host "NOT FOUND: Change1 - formulario"
} else {
  $pr = $rng.Paragraphs(1).Range
  $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="534455F6" w14:textId="464E3E0F" w:rsidR="00857FFC" w:rsidRDefault="00857FFC" w:rsidP="007B30E6"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>El sistema tendrá un formulario el cual tiene que llenar para ir agregando los productos al histórico diario</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> (algunos datos no son obligatorios)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $pr.InsertXML($xml)
}

# Change2 - factura
$rng = $d.Content
$found = $rng.Find.Execute('El dueño podrá compartir la factura por medio de correo o WhatsApp.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  Write-Host "NOT FOUND: Change2 - factura"
} else {
  $pr = $rng.Paragraphs(1).Range
  $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="496DDA17" w14:textId="0D4B3482" w:rsidR="00857FFC" w:rsidRDefault="00857FFC" w:rsidP="007B30E6"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">El dueño podrá compartir la </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>factura</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> por medio de correo o WhatsApp.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $pr.InsertXML($xml)
}

# Change3 - nº de veces
$rng = $d.Content
$found = $rng.Find.Execute('<nº de veces> veces / <unidad de tiempo>', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  Write-Host "NOT FOUND: Change3 - nº de veces"
} else {
  $pr = $rng.Paragraphs(1).Range
  $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="35B6AD26" w14:textId="77777777" w:rsidR="007B30E6" w:rsidRDefault="007B30E6" w:rsidP="007B30E6"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>nº</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> de veces&gt; veces / &lt;unidad de tiempo&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $pr.InsertXML($xml)
}

# Change4a - Urgencia
$rng = $d.Content
$found = $rng.Find.Execute('Urgencia', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  Write-Host "NOT FOUND: Change4a - Urgencia"
} else {
  $pr = $rng.Paragraphs(1).Range
  $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0256D0B0" w14:textId="77777777" w:rsidR="007B30E6" w:rsidRDefault="007B30E6" w:rsidP="007B30E6"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Urgencia</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $pr.InsertXML($xml)
}

# Change4b - Comentarios
$rng = $d.Content
$found = $rng.Find.Execute('Comentarios', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
  Write-Host "NOT FOUND: Change4b - Comentarios"
} else {
  $pr = $rng.Paragraphs(1).Range
  $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="49417FB7" w14:textId="77777777" w:rsidR="007B30E6" w:rsidRDefault="007B30E6" w:rsidP="007B30E6"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/></w:rPr><w:t>Comentarios</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $pr.InsertXML($xml)
}

Write-Host "Done"
